# Update cryptocurrency price/volume data per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.622.08"
$ws.Range("E2").Value = "  +6.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.487.82"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.99"
$ws.Range("E5").Value = "  +7.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.21"
$ws.Range("E6").Value = "  +12.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +7.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.503.80"
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("E10").Value = "  +10.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0985"
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("E12").Value = "  +7.10%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.917.09"
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.525.35"
$ws.Range("E15").Value = "  +5.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.33"
$ws.Range("E16").Value = "  +9.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +6.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.498.92"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("E19").Value = "  +10.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.20"
$ws.Range("E20").Value = "  +10.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.98"
$ws.Range("E21").Value = "  +5.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +10.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.86"
$ws.Range("E24").Value = "  +5.68%  "
$ws.Range("E25").Value = "  +8.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +9.16%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.593.64"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  +9.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0797"
$ws.Range("E30").Value = "  +12.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.96"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +7.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  +5.93%  "
$ws.Range("E36").Value = "  +9.68%  "
$ws.Range("E37").Value = "  +7.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.871"
$ws.Range("E38").Value = "  +11.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.21"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("E40").Value = "  +9.37%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0561"
$ws.Range("E41").Value = "  +8.12%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.615"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.993"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.35"
$ws.Range("E44").Value = "  +10.17%  "
$ws.Range("E45").Value = "  +17.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "260.25"
$ws.Range("E46").Value = "  +20.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0230"
$ws.Range("E47").Value = "  +7.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0920"
$ws.Range("E48").Value = "  +7.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.19"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.903.75"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.69"
$ws.Range("E51").Value = "  +8.10%  "
